$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.536.45'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.13%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.283.62'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.07%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '299.72'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.76%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '95.58'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.86%  '
$ws.Range("E7").Value = '  -2.07%  '
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.490'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.35%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '33.22'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -6.00%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0788'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.33%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '48.79'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.98%  '
$ws.Range("E13").Value = '  +2.01%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '16.55'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +5.84%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.72'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.06%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.641.16'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.83%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.287.05'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.67%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.792'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.50%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '42.464.55'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.08%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0892'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.93%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.42'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.21%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.96'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.40%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '66.72'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.32%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '235.64'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.93%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.95'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.61%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.45'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.58%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '24.21'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.28%  '
$ws.Range("E29").Value = '  -0.56%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '166.70'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.10%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '33.60'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.91%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '9.06'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.73%  '
$ws.Range("E33").Value = '  +0.05%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.69'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.92%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.92'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.78%  '
$ws.Range("B36").Value = 'WEMIXToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.43'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.34%  '
$ws.Range("B37").Value = 'Celestia'
$ws.Range("C37").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '16.92'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.97%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0689'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.66%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.79'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.32%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0995'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.75%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.73'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.03%  '
$ws.Range("E42").Value = '  -2.03%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.38'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.56%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.951.99'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.30%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0278'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.13%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '17.45'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -5.86%  '
$ws.Range("B47").Value = 'FraxShare'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.68'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.81%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.80'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.56%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.511.13'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.64%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '52.39'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -6.17%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.72'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.66%  '
